$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 217; existing rows 217-315 shift down to 219-317.
$ws.Rows("217:218").Insert()

# New row 217: Naranja, Lane Late, Primera
$ws.Range("A217").Value = 5
$ws.Range("B217").Value = "Macroferia Regional de Talca"
$ws.Range("C217").Value = "Maule"
$ws.Range("D217").Value = 44489
$ws.Range("E217").Value = 7
$ws.Range("F217").Value = "Fruta"
$ws.Range("G217").Value = 100102
$ws.Range("H217").Value = "Cítricos"
$ws.Range("I217").Value = 100102005
$ws.Range("J217").Value = "Naranja"
$ws.Range("K217").Value = "Lane Late"
$ws.Range("L217").Value = "Primera"
$ws.Range("M217").Value = 280
$ws.Range("N217").Value = 7000
$ws.Range("O217").Value = 7000
$ws.Range("P217").Value = 7000
$ws.Range("Q217").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R217").Value = "Provincia de Melipilla"
$ws.Range("S217").Value = 467
$ws.Range("T217").Value = 15

# New row 218: Naranja, Navel Late, Primera
$ws.Range("A218").Value = 5
$ws.Range("B218").Value = "Macroferia Regional de Talca"
$ws.Range("C218").Value = "Maule"
$ws.Range("D218").Value = 44489
$ws.Range("E218").Value = 7
$ws.Range("F218").Value = "Fruta"
$ws.Range("G218").Value = 100102
$ws.Range("H218").Value = "Cítricos"
$ws.Range("I218").Value = 100102005
$ws.Range("J218").Value = "Naranja"
$ws.Range("K218").Value = "Navel Late"
$ws.Range("L218").Value = "Primera"
$ws.Range("M218").Value = 460
$ws.Range("N218").Value = 7000
$ws.Range("O218").Value = 7000
$ws.Range("P218").Value = 7000
$ws.Range("Q218").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R218").Value = "Región de O'Higgins"
$ws.Range("S218").Value = 467
$ws.Range("T218").Value = 15
